$wb = $excel.ActiveWorkbook

# --- Sheet "Latest_stocks" (sheet1): update quantities and selection ---
$ws1 = $wb.Worksheets.Item("Latest_stocks")

$ws1.Range("C3").Value = 39
$ws1.Range("C4").Value = 36
$ws1.Range("C6").Value = 4
$ws1.Range("C9").Value = 0
$ws1.Range("C10").Value = 0
$ws1.Range("C12").Value = 0
$ws1.Range("C16").Value = 13
$ws1.Range("C20").Value = 3
$ws1.Range("C32").Value = 11

$ws1.Activate()
$ws1.Range("E20").Select()

# --- Sheet "Arcline_payment" (sheet3): flip status + move selection ---
$ws3 = $wb.Worksheets.Item("Arcline_payment")

$ws3.Range("D17").Value = "PAID"

$ws3.Activate()
$ws3.Range("F23").Select()

$ws1.Activate()
